$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextValue "D2" "29.997.53"
Set-TextValue "E2" "  -0.22%  "
Set-TextValue "D3" "1.897.43"
Set-TextValue "E3" "  -0.85%  "
Set-TextValue "D4" "1.001"
Set-TextValue "E4" "  +0.18%  "
Set-TextValue "D5" "0.8382"
Set-TextValue "E5" "  +3.01%  "
Set-TextValue "D6" "241.86"
Set-TextValue "E6" "  -0.52%  "
Set-TextValue "E7" "  +0.10%  "
Set-TextValue "D8" "0.3295"
Set-TextValue "E8" "  +2.87%  "
Set-TextValue "D9" "26.64"
Set-TextValue "E9" "  +0.89%  "
Set-TextValue "D10" "0.07047"
Set-TextValue "E10" "  +1.25%  "
Set-TextValue "D11" "0.08070"
Set-TextValue "E11" "  +0.39%  "
Set-TextValue "D12" "0.7609"
Set-TextValue "E12" "  +1.09%  "
Set-TextValue "D13" "1.895.86"
Set-TextValue "E13" "  -0.49%  "
Set-TextValue "D14" "5.256"
Set-TextValue "E14" "  +0.17%  "
Set-TextValue "D15" "92.21"
Set-TextValue "E15" "  -1.65%  "
Set-TextValue "D16" "29.991.32"
Set-TextValue "E16" "  -0.18%  "
Set-TextValue "E17" "  +0.30%  "
Set-TextValue "D18" "5.876"
Set-TextValue "E18" "  -2.55%  "
Set-TextValue "D19" "244.18"
Set-TextValue "E19" "  -2.55%  "
Set-TextValue "D20" "0.000007764"
Set-TextValue "E20" "  -0.46%  "
Set-TextValue "D21" "1.000"
Set-TextValue "E21" "  +0.06%  "
Set-TextValue "D22" "2.149.60"
Set-TextValue "E22" "  -0.25%  "
Set-TextValue "D23" "1.001"
Set-TextValue "E23" "  +0.19%  "
Set-TextValue "D24" "6.978"
Set-TextValue "E24" "  -0.15%  "
Set-TextValue "D25" "0.1739"
Set-TextValue "E25" "  +21.79%  "
Set-TextValue "B26" "Monero"
Set-TextValue "C26" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D26" "166.45"
Set-TextValue "E26" "  -1.38%  "
Set-TextValue "B27" "Cosmos"
Set-TextValue "C27" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D27" "9.245"
Set-TextValue "E27" "  -0.95%  "
Set-TextValue "E28" "  -0.78%  "
Set-TextValue "D29" "2.101"
Set-TextValue "E29" "  +1.69%  "
Set-TextValue "D30" "1.361"
Set-TextValue "E30" "  -2.24%  "
Set-TextValue "D31" "1.513"
Set-TextValue "E31" "  -0.85%  "
Set-TextValue "D32" "0.05885"
Set-TextValue "E32" "  +9.64%  "
Set-TextValue "D33" "4.292"
Set-TextValue "E33" "  -2.02%  "
Set-TextValue "D34" "4.077"
Set-TextValue "E34" "  -1.26%  "
Set-TextValue "D35" "1.269"
Set-TextValue "E35" "  +0.04%  "
Set-TextValue "D36" "0.7302"
Set-TextValue "E36" "  -1.61%  "
Set-TextValue "D37" "2.722"
Set-TextValue "E37" "  -0.41%  "
Set-TextValue "D38" "0.01921"
Set-TextValue "E38" "  -0.65%  "
Set-TextValue "D39" "2.779"
Set-TextValue "E39" "  -0.22%  "
Set-TextValue "D40" "0.4439"
Set-TextValue "E40" "  -1.05%  "
Set-TextValue "D41" "72.54"
Set-TextValue "E41" "  -0.72%  "
Set-TextValue "E42" "  -5.43%  "
Set-TextValue "D43" "0.8460"
Set-TextValue "E43" "  +1.55%  "
Set-TextValue "D44" "1.001"
Set-TextValue "E44" "  +0.01%  "
Set-TextValue "E45" "  -1.20%  "
Set-TextValue "D46" "101.64"
Set-TextValue "E46" "  +0.68%  "
Set-TextValue "D47" "1.014.39"
Set-TextValue "E47" "  +5.25%  "
Set-TextValue "B48" "Aptos"
Set-TextValue "C48" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D48" "7.569"
Set-TextValue "E48" "  -1.31%  "
Set-TextValue "B49" "EnergySwap"
Set-TextValue "C49" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D49" "9.809"
Set-TextValue "E49" "  -1.19%  "
Set-TextValue "D50" "2.046.72"
Set-TextValue "E50" "  -0.52%  "
Set-TextValue "D51" "35.91"
Set-TextValue "E51" "  -1.78%  "
